# Update the "Förändrad" (changed) date column (C) for rows 2-10 from
# 45243 (2023-11-13) to 45244 (2023-11-14), as the workbook was refreshed
# by an automatic update process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
